$wb = $excel.ActiveWorkbook
$ts = $wb.TableStyles.Add("MySqlDefault")
Write-Host "created"
